$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '80.952.06'
$ws.Range("E2").Value = '  +2.89%  '
$ws.Range("D3").Value = '3.136.64'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '206.59'
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("D6").Value = '616.44'
$ws.Range("E6").Value = '  -1.97%  '
$ws.Range("D7").Value = '0.278'
$ws.Range("E7").Value = '  +23.31%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.575'
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("D10").Value = '3.136.42'
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").Value = '0.570'
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("D12").Value = '0.0000251'
$ws.Range("E12").Value = '  +13.26%  '
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("D14").Value = '5.23'
$ws.Range("E14").Value = '  -3.17%  '
$ws.Range("D15").Value = '3.718.33'
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("D16").Value = '30.99'
$ws.Range("E16").Value = '  -1.26%  '
$ws.Range("D17").Value = '80.801.83'
$ws.Range("E17").Value = '  +2.71%  '
$ws.Range("D18").Value = '3.141.69'
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("D19").Value = '3.13'
$ws.Range("E19").Value = '  +11.56%  '
$ws.Range("D20").Value = '13.72'
$ws.Range("E20").Value = '  -4.86%  '
$ws.Range("D21").Value = '426.03'
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("D22").Value = '8.85'
$ws.Range("E22").Value = '  -4.89%  '
$ws.Range("D23").Value = '5.02'
$ws.Range("E23").Value = '  +2.29%  '
$ws.Range("D24").Value = '7.16'
$ws.Range("E24").Value = '  +4.89%  '
$ws.Range("D25").Value = '5.10'
$ws.Range("E25").Value = '  +8.16%  '
$ws.Range("D26").Value = '3.307.23'
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").Value = '75.29'
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("D28").Value = '10.69'
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").Value = '0.0000120'
$ws.Range("E30").Value = '  +6.17%  '
$ws.Range("D31").Value = '1.03'
$ws.Range("E31").Value = '  +3.46%  '
$ws.Range("D32").Value = '8.83'
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("D33").Value = '555.13'
$ws.Range("E33").Value = '  +9.51%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = '1.46'
$ws.Range("E34").Value = '  -0.38%  '
$ws.Range("B35").Value = 'Cronos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D35").Value = '0.143'
$ws.Range("E35").Value = '  +11.59%  '
$ws.Range("E36").Value = '  +12.12%  '
$ws.Range("D37").Value = '1.95'
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("D38").Value = '22.47'
$ws.Range("E38").Value = '  -1.28%  '
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("D40").Value = '0.402'
$ws.Range("E40").Value = '  +1.71%  '
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").Value = '20.69'
$ws.Range("E41").Value = '  +3.75%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").Value = '5.86'
$ws.Range("E42").Value = '  +9.14%  '
$ws.Range("D43").Value = '2.98'
$ws.Range("E43").Value = '  +20.67%  '
$ws.Range("D44").Value = '1.98'
$ws.Range("E44").Value = '  +11.93%  '
$ws.Range("D45").Value = '159.33'
$ws.Range("E45").Value = '  -2.83%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").Value = '185.61'
$ws.Range("E47").Value = '  -2.86%  '
$ws.Range("D48").Value = '44.93'
$ws.Range("E48").Value = '  +5.55%  '
$ws.Range("D49").Value = '1.30'
$ws.Range("E49").Value = '  +1.19%  '
$ws.Range("D50").Value = '0.759'
$ws.Range("E50").Value = '  -5.78%  '
$ws.Range("D51").Value = '25.29'
$ws.Range("E51").Value = '  +2.91%  '
